# Round the ConvexHullArea values in column D (rows 2-26) to the nearest
# whole number, as plain numeric values (not just display formatting).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

for ($row = 2; $row -le 26; $row++) {
    $cell = $ws.Cells.Item($row, 4)
    $current = [double]$cell.Value()
    $cell.Value = $excel.WorksheetFunction.Round($current, 0)
}
